$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (column F) counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 14028
$wsExhibit.Range("F5").Value = 247
$wsExhibit.Range("F9").Value = 58
$wsExhibit.Range("F10").Value = 567
$wsExhibit.Range("F14").Value = 14342
$wsExhibit.Range("F15").Value = 386
$wsExhibit.Range("F16").Value = 644
$wsExhibit.Range("F17").Value = 15086
$wsExhibit.Range("F19").Value = 8468
$wsExhibit.Range("F22").Value = 45
$wsExhibit.Range("F25").Value = 173
$wsExhibit.Range("F27").Value = 16
$wsExhibit.Range("F29").Value = 44
$wsExhibit.Range("F31").Value = 2
$wsExhibit.Range("F32").Value = 2
$wsExhibit.Range("F33").Value = 38
$wsExhibit.Range("F34").Value = 22
$wsExhibit.Range("F35").Value = 27
$wsExhibit.Range("F36").Value = 409
$wsExhibit.Range("F37").Value = 14
$wsExhibit.Range("F38").Value = 18
$wsExhibit.Range("F40").Value = 241
$wsExhibit.Range("F41").Value = 402
$wsExhibit.Range("F43").Value = 5214

# Sheet "全部类型" - update "想去人数" (column F) counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 14028
$wsAll.Range("F5").Value = 247
$wsAll.Range("F9").Value = 58
$wsAll.Range("F10").Value = 567
$wsAll.Range("F14").Value = 14342
$wsAll.Range("F15").Value = 386
$wsAll.Range("F16").Value = 644
$wsAll.Range("F17").Value = 15086
$wsAll.Range("F19").Value = 8468
$wsAll.Range("F23").Value = 45
$wsAll.Range("F26").Value = 173
$wsAll.Range("F28").Value = 16
$wsAll.Range("F30").Value = 44
$wsAll.Range("F32").Value = 2
$wsAll.Range("F33").Value = 2
$wsAll.Range("F34").Value = 38
$wsAll.Range("F35").Value = 22
$wsAll.Range("F36").Value = 27
$wsAll.Range("F39").Value = 409
$wsAll.Range("F40").Value = 14
$wsAll.Range("F41").Value = 18
$wsAll.Range("F43").Value = 241
$wsAll.Range("F44").Value = 402
$wsAll.Range("F46").Value = 5214
